# Update crypto price/volume snapshot figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "307.63" },
    @{ Cell = "E2";  Value = "-0.63%" },
    @{ Cell = "D3";  Value = "37.14" },
    @{ Cell = "E3";  Value = "-0.06%" },
    @{ Cell = "D4";  Value = "5.115" },
    @{ Cell = "E4";  Value = "-0.23%" },
    @{ Cell = "D5";  Value = "0.07828" },
    @{ Cell = "E5";  Value = "0.77%" },
    @{ Cell = "D6";  Value = "4.399" },
    @{ Cell = "E6";  Value = "0.13%" },
    @{ Cell = "D7";  Value = "8.246" },
    @{ Cell = "E7";  Value = "0.49%" },
    @{ Cell = "D8";  Value = "1.886" },
    @{ Cell = "E8";  Value = "0.51%" },
    @{ Cell = "D9";  Value = "2.999" },
    @{ Cell = "E9";  Value = "-1.62%" },
    @{ Cell = "D10"; Value = "0.9207" },
    @{ Cell = "E10"; Value = "0.30%" },
    @{ Cell = "E11"; Value = "-9.00%" },
    @{ Cell = "D12"; Value = "0.1895" },
    @{ Cell = "E12"; Value = "0.29%" },
    @{ Cell = "D13"; Value = "0.08855" },
    @{ Cell = "E13"; Value = "-4.43%" },
    @{ Cell = "D14"; Value = "0.03350" },
    @{ Cell = "E14"; Value = "-2.56%" },
    @{ Cell = "D15"; Value = "0.09578" },
    @{ Cell = "E15"; Value = "-1.02%" },
    @{ Cell = "D16"; Value = "0.001378" },
    @{ Cell = "E16"; Value = "0.47%" },
    @{ Cell = "D17"; Value = "0.005684" },
    @{ Cell = "E17"; Value = "-0.83%" },
    @{ Cell = "D18"; Value = "3.409" },
    @{ Cell = "E18"; Value = "-4.07%" },
    @{ Cell = "D19"; Value = "0.3423" },
    @{ Cell = "D20"; Value = "6.291" },
    @{ Cell = "E20"; Value = "19.36%" },
    @{ Cell = "D21"; Value = "0.1286" },
    @{ Cell = "E21"; Value = "0.87%" },
    @{ Cell = "D22"; Value = "0.2452" },
    @{ Cell = "E22"; Value = "-5.45%" },
    @{ Cell = "D23"; Value = "0.04344" },
    @{ Cell = "E23"; Value = "0.16%" },
    @{ Cell = "D24"; Value = "0.001193" },
    @{ Cell = "E24"; Value = "-0.48%" },
    @{ Cell = "D25"; Value = "0.004262" },
    @{ Cell = "E25"; Value = "0.40%" },
    @{ Cell = "E26"; Value = "6.81%" },
    @{ Cell = "E27"; Value = "-98.10%" },
    @{ Cell = "D39"; Value = "0.02170" },
    @{ Cell = "E39"; Value = "4.99%" },
    @{ Cell = "D40"; Value = "0.05018" },
    @{ Cell = "E40"; Value = "-0.56%" },
    @{ Cell = "D41"; Value = "0.007549" },
    @{ Cell = "E41"; Value = "-1.73%" },
    @{ Cell = "D42"; Value = "0.1349" },
    @{ Cell = "E42"; Value = "0.36%" },
    @{ Cell = "D43"; Value = "0.008671" },
    @{ Cell = "E43"; Value = "-11.81%" },
    @{ Cell = "D44"; Value = "0.002040" },
    @{ Cell = "E44"; Value = "-6.12%" },
    @{ Cell = "D45"; Value = "0.008731" },
    @{ Cell = "E45"; Value = "-0.36%" },
    @{ Cell = "D46"; Value = "0.00006534" },
    @{ Cell = "E46"; Value = "-2.63%" },
    @{ Cell = "E47"; Value = "-0.05%" },
    @{ Cell = "D48"; Value = "0.003296" },
    @{ Cell = "E48"; Value = "12.20%" },
    @{ Cell = "E49"; Value = "-16.58%" },
    @{ Cell = "D50"; Value = "0.00002101" },
    @{ Cell = "E50"; Value = "-0.05%" },
    @{ Cell = "D51"; Value = "0.0002001" },
    @{ Cell = "E51"; Value = "-0.05%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
